$d = $word.ActiveDocument

# The first table on the page holds the "Standard Solution #1/#2/#3" rows.
# Column 3 ("Absorbance at ___ nm") is currently empty (highlighted yellow)
# for rows 4-6 and needs the measured absorbance values typed in, using the
# same run formatting (theme fonts + szCs) already used throughout the
# table's paragraph marks.

function Fill-AbsorbanceCell {
    param(
        [int]$RowIndex,
        [string]$ParaId,
        [string]$Value
    )

    $table = $d.Tables.Item(1)
    $cell = $table.Cell($RowIndex, 3)
    $rng = $cell.Range
    # Collapse to an insertion point right before the cell's end-of-paragraph
    # mark (the cell is currently empty, so Start/End bracket just that mark).
    $rng.End = $rng.End - 1

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
        '<w:body>' + `
        '<w:p w14:paraId="' + $ParaId + '" w14:textId="77777777" w:rsidR="0028599E" w:rsidRPr="007C71A3" w:rsidRDefault="0028599E" w:rsidP="00801E5E">' + `
        '<w:pPr><w:spacing w:after="0" w:line="276" w:lineRule="auto"/><w:jc w:val="left"/>' + `
        '<w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:szCs w:val="24"/></w:rPr>' + `
        '</w:pPr>' + `
        '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:szCs w:val="24"/></w:rPr><w:t>' + $Value + '</w:t></w:r>' + `
        '</w:p>' + `
        '</w:body></w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'

    $rng.InsertXML($xml)
}

# Note: this runtime's PowerShell interpreter does not bind named
# (-Param value) arguments correctly, so call positionally.

# Row 4 = "Standard Solution #1" -> Absorbance 0.152
Fill-AbsorbanceCell 4 "68C43957" "0.152"

# Row 5 = "Standard Solution #2" -> Absorbance 0.737
Fill-AbsorbanceCell 5 "6EC29BD6" "0.737"

# Row 6 = "Standard Solution #3" -> Absorbance 1.363
Fill-AbsorbanceCell 6 "3137359B" "1.363"

Write-Host "done"
